# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages have now been generated: status moves from
# "Ready for handoff" to "Handed back: in sync with en-US", the "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns get populated, and a couple of columns get widened so the new
# long file names are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This touches the Overview summary columns (E/F) for both rows, and
#    the per-language "Status" column (C) on the zh-cn / de-de sheets.
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Column widths: widen the columns that now hold long handback file
#    names / datetimes so they are fully visible.
# ---------------------------------------------------------------------
# Overview: zh-cn (E) / de-de (F) status columns
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de: Status (C), Latest Target File (I), Latest Handback File (J)
foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(3).ColumnWidth  = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth  = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

# ---------------------------------------------------------------------
# 3. Populate the handback results for each language / each row.
#    Column I = Latest Target File, J = Latest Handback File,
#    K = Latest Handback DateTime.
# ---------------------------------------------------------------------
$sourceMd     = "6808ca51-4609-41ca-bd42-6a75a51e4a3e.md"
$sourceMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5e2833389adc6f1a57b5d3d1c673febbb3d40b/e2e/6808ca51-4609-41ca-bd42-6a75a51e4a3e.md"

$zhXlf  = "6808ca51-4609-41ca-bd42-6a75a51e4a3e.ae41ee51ac863b8283b38ce29343d0bb75a8253e.zh-cn.xlf"
$deXlf  = "6808ca51-4609-41ca-bd42-6a75a51e4a3e.ae41ee51ac863b8283b38ce29343d0bb75a8253e.de-de.xlf"

$zhHandbackTime = "2016-09-06 11:23:32"
$deHandbackTime = "2016-09-06 11:23:41"

# -- zh-cn sheet, rows 2 and 3 --
$zhcn.Range("I2").Value = $sourceMd
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $sourceMdUrl, "", "", $sourceMd) | Out-Null
$zhcn.Range("J2").Value = $zhXlf
$zhcn.Range("K2").Value = $zhHandbackTime

$zhcn.Range("I3").Value = $sourceMd
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $sourceMdUrl, "", "", $sourceMd) | Out-Null
$zhcn.Range("J3").Value = $zhXlf
$zhcn.Range("K3").Value = $zhHandbackTime

# -- de-de sheet, rows 2 and 3 --
$dede.Range("I2").Value = $sourceMd
$dede.Hyperlinks.Add($dede.Range("I2"), $sourceMdUrl, "", "", $sourceMd) | Out-Null
$dede.Range("J2").Value = $deXlf
$dede.Range("K2").Value = $deHandbackTime

$dede.Range("I3").Value = $sourceMd
$dede.Hyperlinks.Add($dede.Range("I3"), $sourceMdUrl, "", "", $sourceMd) | Out-Null
$dede.Range("J3").Value = $deXlf
$dede.Range("K3").Value = $deHandbackTime
